$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price values that look numeric must be forced to stay as text
# (matching the original inlineStr/General-format cells) by prefixing with
# a literal apostrophe (Excel quote-prefix) and then resetting the cell style
# back to Normal/General so no quotePrefix style flag lingers.

$ws.Range("D2").Value = '69.594.10'
$ws.Range("E2").Value = '  -0.01%  '
$ws.Range("D3").Value = '2.507.02'
$ws.Range("E3").Value = '  -0.08%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'575.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.30%  '
$ws.Range("D6").Value = "'166.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = "'0.514"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.24%  '
$ws.Range("D9").Value = '2.506.25'
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("D10").Value = "'0.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.75%  '
$ws.Range("E11").Value = '  -0.39%  '
$ws.Range("D12").Value = "'0.356"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.88%  '
$ws.Range("D13").Value = "'4.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.04%  '
$ws.Range("D14").Value = '2.965.46'
$ws.Range("E14").Value = '  -0.58%  '
$ws.Range("D15").Value = '69.472.38'
$ws.Range("E15").Value = '  +0.07%  '
$ws.Range("D16").Value = "'0.0000178"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.30%  '
$ws.Range("D17").Value = "'24.86"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.44%  '
$ws.Range("D18").Value = '2.504.75'
$ws.Range("E18").Value = '  -0.44%  '
$ws.Range("D19").Value = "'11.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.64%  '
$ws.Range("D20").Value = "'7.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.22%  '
$ws.Range("D21").Value = "'349.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.20%  '
$ws.Range("E22").Value = '  -1.18%  '
$ws.Range("E23").Value = '  -1.21%  '
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("D25").Value = "'70.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.79%  '
$ws.Range("D26").Value = "'3.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.68%  '
$ws.Range("D27").Value = "'8.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.14%  '
$ws.Range("D28").Value = '2.640.57'
$ws.Range("E28").Value = '  +0.07%  '
$ws.Range("D29").Value = "'0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("D30").Value = '0.0₃0889'
$ws.Range("E30").Value = '  -1.54%  '
$ws.Range("E31").Value = '  -1.08%  '
$ws.Range("D32").Value = "'461.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.74%  '
$ws.Range("E33").Value = '  -5.66%  '
$ws.Range("E34").Value = '  -1.22%  '
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("D36").Value = "'159.90"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.17%  '
$ws.Range("E37").Value = '  +0.12%  '
$ws.Range("D38").Value = "'19.06"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.81%  '
$ws.Range("D39").Value = "'18.47"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.45%  '
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("D41").Value = "'0.319"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.27%  '
$ws.Range("D42").Value = "'4.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.75%  '
$ws.Range("D43").Value = "'1.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.96%  '
$ws.Range("D44").Value = "'38.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("E45").Value = '  -5.13%  '
$ws.Range("E46").Value = '  -8.03%  '
$ws.Range("D47").Value = "'142.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.66%  '
$ws.Range("D48").Value = "'3.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.20%  '
$ws.Range("E49").Value = '  -2.44%  '
$ws.Range("E50").Value = '  +0.48%  '
$ws.Range("D51").Value = "'0.580"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.26%  '
